$d = $word.ActiveDocument

# Locate the paragraph that ends the "General Resources" Linux-as-desktop
# blurb ("... Linux as a regular operating system on your machine.") so we
# can insert the two new "Linux as a router (Fedora)" paragraphs right
# after it, and drop the paragraph's own "space after" (w:after="0").
$anchor = $d.Content
$anchor.Find.Execute("Linux as a regular operating system on your machine.", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$para = $anchor.Paragraphs(1)
$target = $para.Range
$target.Collapse(0)

$ns = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

$newXml = '<w:p xmlns:w="' + $ns + '">' + `
            '<w:pPr><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr>' + `
            '<w:r><w:t>Tangentially related to the project but contains some good information if you want to use Linux as a regular operating system on your machine.</w:t></w:r>' + `
          '</w:p>' + `
          '<w:p xmlns:w="' + $ns + '">' + `
            '<w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/></w:pPr>' + `
            '<w:r><w:t xml:space="preserve">Linux as a router (Fedora): </w:t></w:r>' + `
            '<w:r><w:t>FEDORA_ROUTER_LINK_PLACEHOLDER</w:t></w:r>' + `
            '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
          '</w:p>' + `
          '<w:p xmlns:w="' + $ns + '">' + `
            '<w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/></w:pPr>' + `
            '<w:r><w:tab/><w:t>A guide on configuring Fedora Server as a router gateway.</w:t></w:r>' + `
          '</w:p>'

$target.InsertXML($newXml)

# Turn the placeholder run into a real hyperlink (this mints the
# relationship + applies the built-in "Hyperlink" character style, just
# like Word does when you paste/insert a link).
$linkRange = $d.Content
$linkRange.Find.Execute("FEDORA_ROUTER_LINK_PLACEHOLDER", $true, $false, $false, `
    $false, $false, $true, 1, $false, "", 0) | Out-Null

$d.Hyperlinks.Add($linkRange, `
    "https://fedoramagazine.org/use-fedora-server-create-router-gateway/", `
    $null, $null, `
    "https://fedoramagazine.org/use-fedora-server-create-router-gateway/") | Out-Null

Write-Output "Inserted Fedora router paragraphs"
